# FatsGuide.pptx edit:
#  1. Add alternate-text (alt text) descriptions to the nutrition-label /
#     food-tracker pictures on every food slide.
#  2. Remove the duplicate "Almonds" slide (position 3).
#  3. Refresh the cached "datetimeFigureOut" footer date field
#     (3/29/2020 -> 4/6/2020) on the slide master and every slide layout.

$p = $ppt.ActivePresentation

# --- 1a. Almonds (slide 2) -------------------------------------------------
$s = $p.Slides.Item(2)
$s.Shapes.Item(1).AlternativeText = "Almonds food tracker showing 75% fat, 11% carbs, and 14% protein"
$s.Shapes.Item(2).AlternativeText = "Almonds nutritional label: Serving size 1 oz., calories 170, total fat 15g (26%), saturated fat 1g, polyunsaturated fat 4g, monounsaturated fat 10g, total carbohydrates 5g (5%), dietary fiber 3g, sugars 1g, proteins 6g (4%)"

# --- 2. Delete the duplicate Almonds slide (slide 3) -----------------------
$p.Slides.Item(3).Delete()

# After the delete, the remaining food slides (formerly 4-10) shift up to
# become slides 3-9, in the same relative order.

# --- 1b. Avocado (now slide 3) ---------------------------------------------
$s = $p.Slides.Item(3)
$s.Shapes.Item(1).AlternativeText = "Avocado nutritional label: Serving size 1/2 medium, calories 120, total fat 11g (19%), saturated fat 1g, polyunsaturated fat 4g, monounsaturated fat 10g, total carbohydrates 6g (7%), dietary fiber 3g, sugars 1g, protein 1.5g (1%)"
$s.Shapes.Item(2).AlternativeText = "Avocado food tracker showing 76% fat, 20% carbs, and 4% protein"

# --- 1c. Bacon (now slide 4) ------------------------------------------------
$s = $p.Slides.Item(4)
$s.Shapes.Item(1).AlternativeText = "Bacon nutritional label: Serving size 1 piece, calories 60, total fat 5g (9%), saturated fat 2g, polyunsaturated fat 0g, monounsaturated fat 0g, total carbohydrates 0g (0%), dietary fiber 0g, sugars 0g, proteins 3g (2%)"
$s.Shapes.Item(2).AlternativeText = "Bacon food tracker showing 79% fat, 0% carbs, and 21% protein"

# --- 1d. Butter (now slide 5) -----------------------------------------------
$s = $p.Slides.Item(5)
$s.Shapes.Item(1).AlternativeText = "Butter nutritional label: Serving size 1 oz., calories 120, total fat 9g (16%), saturated fat 6 g, polyunsaturated fat 0g, monounsaturated fat 0g, total carbohydrates 2g (2%), dietary fiber 0g, sugars 0g, proteins 6g (4%)"
$s.Shapes.Item(2).AlternativeText = "Butter food tracker showing 100% fat, 0% carbs, and 0% protein`n"

# --- 1e. Cashews (now slide 6) ----------------------------------------------
$s = $p.Slides.Item(6)
$s.Shapes.Item(1).AlternativeText = "Casheews nutritional label: Serving size 1 oz., calories 170, total fat 14g (24%), saturated fat 2g, polyunsaturated fat 0g, monounsaturated fat 0g, total carbohydrates 6g (6%), dietary fiber 2g, sugars 1g, proteins 8g (5%)"
$s.Shapes.Item(2).AlternativeText = "Cashews food tracker showing 69% fat, 13% carbs, and 18% protein"

# --- 1f. Cheddar cheese (now slide 7) ---------------------------------------
$s = $p.Slides.Item(7)
$s.Shapes.Item(1).AlternativeText = "Cheddar cheese: Serving size 1 oz., calories 120, total fat 9g (16%), saturated fat 6 g, polyunsaturated fat 0g, monounsaturated fat 0g, total carbohydrates 2g (2%), dietary fiber 0g, sugars 0g, proteins 6g (4%)"
$s.Shapes.Item(2).AlternativeText = "Cheddar cheese food tracker showing 72% fat, 7% carbs, and 21% protein"

# --- 1g. Egg (now slide 8) --------------------------------------------------
$s = $p.Slides.Item(8)
$s.Shapes.Item(1).AlternativeText = "Egg nutritional label: Serving size 1 large., calories 70, total fat 5g (8%), saturated fat 1.5g, polyunsaturated fat 1g, monounsaturated fat 2g, total carbohydrates 0.5g (0%), dietary fiber 0g, sugars 0g, proteins 6g (4%)"
$s.Shapes.Item(2).AlternativeText = "Egg food tracker showing 62% fat, 2% carbs, and 36% protein`n"

# --- 1h. Peanut butter (now slide 9) ----------------------------------------
$s = $p.Slides.Item(9)
$s.Shapes.Item(1).AlternativeText = "Peanut butter nutritional label: Serving size 2 Tbsp., calories 190, total fat 16g (28%), saturated fat 2.5g, polyunsaturated fat 0g, monounsaturated fat 0g, total carbohydrates 8g (9%), dietary fiber 2g, sugars 3g, proteins 7g (5%)"
$s.Shapes.Item(2).AlternativeText = "Peanut butter food tracker showing 71% fat, 16% carbs, and 13% protein"

# --- 3. Refresh the cached footer date field --------------------------------
$m = $p.SlideMaster
for ($j = 1; $j -le $m.Shapes.Count; $j++) {
  $sh = $m.Shapes.Item($j)
  if ($sh.Name -like "Date Placeholder*") {
    $sh.TextFrame.TextRange.Text = "4/6/2020"
  }
}
for ($i = 1; $i -le $m.CustomLayouts.Count; $i++) {
  $cl = $m.CustomLayouts.Item($i)
  for ($j = 1; $j -le $cl.Shapes.Count; $j++) {
    $sh = $cl.Shapes.Item($j)
    if ($sh.Name -like "Date Placeholder*") {
      $sh.TextFrame.TextRange.Text = "4/6/2020"
    }
  }
}
